# commit thêm chức năng cập nhật khấu trừ....
# Update "Lương tăng ca" (col D) and "Tổng lương nhân viên" (col F) values
# for several employees to reflect updated deduction calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("D3").Value = 1500

$ws.Range("D4").Value = 0
$ws.Range("F4").Value = 2000

$ws.Range("D7").Value = 2000
$ws.Range("F7").Value = 2000

$ws.Range("F8").Value = 1500

$ws.Range("D9").Value = 1500
$ws.Range("F9").Value = 2000
